# Added K Means Demo
# Append 5 new sample rows (Opportunity_ID 1021-1025) to the SalesOpportunityDataSet
# sheet, mirroring the existing data layout, and wrap/vertically-center the
# new Contact_Title cells so the longer job titles read cleanly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(1021, "Manufacturing", "Large",  "Operations Director",   92, "IoT Monitoring",  "West",        "Yes", 1.8,  1),
    @(1022, "Education",     "Medium", "Technology Specialist", 51, "LMS Integration", "Northeast",   "No",  9,    0),
    @(1023, "Government",    "Small",  "Cybersecurity Officer", 70, "Secure Cloud",    "Midwest",     "No",  6,    0),
    @(1024, "Finance",       "Large",  "Director of FinOps",    85, "Risk Analytics",  "Mid-Atlantic","Yes", 2.1,  1),
    @(1025, "Retail",        "Medium", "E-commerce Manager",    45, "CRM Suite",       "Southeast",   "No",  11.5, 0)
)

$startRow = 22
$tallRows = @(23, 24, 26)   # rows whose Contact_Title text needs two display lines

$r = $startRow
foreach ($row in $rows) {
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}

$lastRow = $r - 1
$dataRange = $ws.Range("A$($startRow):J$lastRow")

# Build the wrap-text / vertical-center format on the first cell only, then
# fan it out to the rest of the new block with a format-only paste so the
# whole range ends up sharing a single new cellXfs entry.
$firstCell = $ws.Range("A$startRow")
$firstCell.WrapText = $true
$firstCell.VerticalAlignment = -4108   # xlCenter
$firstCell.Copy()
$dataRange.PasteSpecial(-4122)          # xlPasteFormats

foreach ($tr in $tallRows) {
    $ws.Rows.Item($tr).RowHeight = 30
}

$ws.Range("I25").Select()
